$d = $word.ActiveDocument

# --- Edit 1: merge the three runs around the spell-checked word "sul" into
# a single run (removing the now-unnecessary w:proofErr markers). Performing
# a Find/Replace over the full phrase (including the smart-quoted "sul")
# causes Word to re-emit the matched text as one run.
$quote = [char]0x201C
$rquote = [char]0x201D
$oldPhrase = "We noticed both " + $quote + "sul" + $rquote + " and roman numerals"
$newPhrase = "We noticed both " + $quote + "sul" + $rquote + " and roman numerals"
$d.Content.Find.Execute($oldPhrase, $true, $false, $false, $false, $false, $true, 1, $false, $newPhrase, 2) | Out-Null

# --- Edit 2: append four new errata rows (mm. 455, 469, 473, 480) to the
# last table, cloning the formatting of the table's current last row.
$t = $d.Tables.Item(1)

function Add-ErrataRow($instrument, $measure, $note, $heightPts) {
    $row = $t.Rows.Add()
    $row.Height = $heightPts
    $row.Cells.Item(1).Range.Text = $instrument
    $row.Cells.Item(2).Range.Text = $measure
    $row.Cells.Item(3).Range.Text = $note
    # Leave the last ("resolved?") cell empty - typing then deleting a
    # placeholder character clears the auto-cloned empty run so the cell
    # matches the rest of the table's blank cells exactly.
    $lastCell = $row.Cells.Item(4)
    $lastCell.Range.Text = "X"
    $lastCell.Range.Characters.Item(1).Delete()
}

Add-ErrataRow "Violin 1" "455" "At which point would you like the violin to return to arco? " 49.2
Add-ErrataRow "Violin 1" "469" "Are the two consecutive G naturals in the grace note group intentional? If not, what should the new note be?" 49.2
Add-ErrataRow "Violin 1" "473" "Would you like the slur to cover all the notes in this scale, or just the first four?" 42.1
Add-ErrataRow "Violin 2" "480" "Semiquavers have been regrouped for ease of reading" 42.1
